$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H29").Value = 2375
$ws.Range("I29").Value = 500
$ws.Range("J29").Value = 3000
$ws.Range("K29").Value = 1500
$ws.Range("L29").Value = 9000
$ws.Range("M29").Value = -1219
$ws.Range("N29").Value = -9562
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").ClearContents()
$ws.Range("H43").Value = 17597.4
$ws.Range("J43").Value = 18441.555
$ws.Range("L43").Value = 18441.555
$ws.Range("N43").Value = -18579.555
$ws.Range("H46").Value = 3609.375
$ws.Range("H60").Value = 3609.375
$ws.Range("H70").Value = 2483.111
$ws.Range("I70").Value = 2387.5
$ws.Range("J70").Value = 2559.6
$ws.Range("K70").Value = 7162.5
$ws.Range("L70").Value = 7678.799999999999
$ws.Range("M70").Value = -6892.5
$ws.Range("N70").Value = -8218.799999999999
$ws.Range("H73").Value = 2483.111
$ws.Range("I73").Value = 2387.5
$ws.Range("J73").Value = 2559.6
$ws.Range("K73").Value = 7162.5
$ws.Range("L73").Value = 7678.799999999999
$ws.Range("M73").Value = -6226.5
$ws.Range("N73").Value = -9550.799999999999
$ws.Range("H96").Value = 1247
$ws.Range("J96").Value = 1383.6
$ws.Range("L96").Value = 4150.799999999999
$ws.Range("N96").Value = -6896.799999999999
$ws.Range("H100").Value = 1644.5769
$ws.Range("J100").Value = 2070.25
$ws.Range("L100").Value = 2070.25
$ws.Range("N100").Value = -3152.25
$ws.Range("H107").Value = 778.5454999999999
$ws.Range("I107").Value = 451.625
$ws.Range("K107").Value = 451.625
$ws.Range("M107").Value = 1468.375
$ws.Range("H116").Value = 2656.681
$ws.Range("I116").Value = 2408.7815
$ws.Range("J116").Value = 5737.7144
$ws.Range("K116").Value = 2408.7815
$ws.Range("L116").Value = 5737.7144
$ws.Range("M116").Value = 1033.2185
$ws.Range("N116").Value = -12621.7144
$ws.Range("H132").Value = 6572.576
$ws.Range("I132").Value = 4859.582
$ws.Range("K132").Value = 14578.746
$ws.Range("M132").Value = -12048.746
$ws.Range("H135").Value = 7238.6523
$ws.Range("I135").Value = 7802.4443
$ws.Range("J135").Value = 6876.2144
$ws.Range("K135").Value = 70221.9987
$ws.Range("L135").Value = 61885.9296
$ws.Range("M135").Value = -67686.9987
$ws.Range("N135").Value = -66955.9296
$ws.Range("H137").Value = 10880.568
$ws.Range("I137").Value = 3355
$ws.Range("J137").Value = 22832.941
$ws.Range("K137").Value = 10065
$ws.Range("L137").Value = 68498.823
$ws.Range("M137").Value = -7515
$ws.Range("N137").Value = -73598.823
$ws.Range("H139").Value = 85856.86
$ws.Range("J139").Value = 87666.336
$ws.Range("L139").Value = 87666.336
$ws.Range("N139").Value = -97946.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 814.2778
$ws.Range("I2").Value = 490.25
$ws.Range("J2").Value = 1462.3334
$ws.Range("K2").Value = 490.25
$ws.Range("L2").Value = 1462.3334
$ws.Range("M2").Value = -377.25
$ws.Range("N2").Value = -1688.3334
$ws.Range("H74").Value = 13218.811
$ws.Range("I74").Value = 2217.5454
$ws.Range("K74").Value = 2217.5454
$ws.Range("M74").Value = -1343.5454
$ws.Range("H77").Value = 13218.811
$ws.Range("I77").Value = 2217.5454
$ws.Range("K77").Value = 11087.727
$ws.Range("M77").Value = -6719.726999999999
$ws.Range("H116").Value = 814.2778
$ws.Range("I116").Value = 490.25
$ws.Range("J116").Value = 1462.3334
$ws.Range("K116").Value = 490.25
$ws.Range("L116").Value = 1462.3334
$ws.Range("M116").Value = 1803.75
$ws.Range("N116").Value = -6050.3334
$ws.Range("H132").Value = 2008069.5
$ws.Range("I132").Value = 3309.7568
$ws.Range("J132").Value = 7713924.5
$ws.Range("K132").Value = 9929.270400000001
$ws.Range("L132").Value = 23141773.5
$ws.Range("M132").Value = -7399.270400000001
$ws.Range("N132").Value = -23146833.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 814.2778
$ws.Range("I3").Value = 490.25
$ws.Range("J3").Value = 1462.3334
$ws.Range("K3").Value = 490.25
$ws.Range("L3").Value = 1462.3334
$ws.Range("M3").Value = -376.25
$ws.Range("N3").Value = -1690.3334
$ws.Range("H99").Value = 3434.111
$ws.Range("I99").Value = 3570.8572
$ws.Range("J99").Value = 2955.5
$ws.Range("K99").Value = 3570.8572
$ws.Range("L99").Value = 2955.5
$ws.Range("M99").Value = -2072.8572
$ws.Range("N99").Value = -5951.5
$ws.Range("H105").Value = 5292.5713
$ws.Range("I105").Value = 5292.5713
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 5292.5713
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -3545.5713
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 2529.9119
$ws.Range("I107").Value = 2751.1072
$ws.Range("K107").Value = 2751.1072
$ws.Range("M107").Value = -831.1071999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16719
$ws.Range("I31").Value = 5785.423
$ws.Range("K31").Value = 5785.423
$ws.Range("M31").Value = -5490.423
$ws.Range("H34").Value = 16719
$ws.Range("I34").Value = 5785.423
$ws.Range("K34").Value = 5785.423
$ws.Range("M34").Value = -5583.423
$ws.Range("H107").Value = 901.2121
$ws.Range("I107").Value = 735.28
$ws.Range("K107").Value = 735.28
$ws.Range("M107").Value = 1184.72
$ws.Range("H122").Value = 6233.4165
$ws.Range("J122").Value = 9447
$ws.Range("L122").Value = 28341
$ws.Range("N122").Value = -33241
$ws.Range("H132").Value = 7583.483
$ws.Range("I132").Value = 3006.3157
$ws.Range("J132").Value = 16280.1
$ws.Range("K132").Value = 9018.947100000001
$ws.Range("L132").Value = 48840.3
$ws.Range("M132").Value = -6488.947100000001
$ws.Range("N132").Value = -53900.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 2837.25
$ws.Range("I29").Value = 2837.25
$ws.Range("K29").Value = 8511.75
$ws.Range("M29").Value = -8234.75
$ws.Range("H80").Value = 18324.9
$ws.Range("J80").Value = 18324.9
$ws.Range("L80").Value = 54974.7
$ws.Range("N80").Value = -56846.7
$ws.Range("H83").Value = 18324.9
$ws.Range("J83").Value = 18324.9
$ws.Range("L83").Value = 164924.1
$ws.Range("N83").Value = -174284.1
$ws.Range("H108").Value = 4644.25
$ws.Range("I108").Value = 384.66666
$ws.Range("J108").Value = 7200
$ws.Range("K108").Value = 1153.99998
$ws.Range("L108").Value = 21600
$ws.Range("M108").Value = 1726.00002
$ws.Range("N108").Value = -27360
$ws.Range("H109").Value = 2391440
$ws.Range("I109").Value = 5118.4287
$ws.Range("J109").Value = 4777761.5
$ws.Range("K109").Value = 15355.2861
$ws.Range("L109").Value = 14333284.5
$ws.Range("M109").Value = -14315.2861
$ws.Range("N109").Value = -14335364.5
$ws.Range("H112").Value = 16175
$ws.Range("J112").Value = 16763.158
$ws.Range("L112").Value = 50289.474
$ws.Range("N112").Value = -52505.474
$ws.Range("H115").Value = 4760.8887
$ws.Range("I115").Value = 771.6
$ws.Range("J115").Value = 9747.5
$ws.Range("K115").Value = 2314.8
$ws.Range("L115").Value = 29242.5
$ws.Range("M115").Value = -1139.8
$ws.Range("N115").Value = -31592.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11674.654
$ws.Range("I80").Value = 8394.1875
$ws.Range("K80").Value = 8394.1875
$ws.Range("M80").Value = -7396.1875
$ws.Range("H83").Value = 11674.654
$ws.Range("I83").Value = 8394.1875
$ws.Range("K83").Value = 41970.9375
$ws.Range("M83").Value = -36978.9375
$ws.Range("H102").Value = 2617.6
$ws.Range("I102").Value = 2686.2222
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2686.2222
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -1064.2222
$ws.Range("N102").Value = -5244
$ws.Range("H122").Value = 1875.091
$ws.Range("J122").Value = 3073
$ws.Range("L122").Value = 9219
$ws.Range("N122").Value = -14119
$ws.Range("H126").Value = 6127.3335
$ws.Range("I126").Value = 5001.375
$ws.Range("J126").Value = 7414.143
$ws.Range("K126").Value = 15004.125
$ws.Range("L126").Value = 22242.429
$ws.Range("M126").Value = -12534.125
$ws.Range("N126").Value = -27182.429
$ws.Range("H127").Value = 83333
$ws.Range("J127").Value = 83333
$ws.Range("L127").Value = 83333
$ws.Range("N127").Value = -93253
$ws.Range("H132").Value = 49173.04
$ws.Range("I132").Value = 59214.26
$ws.Range("K132").Value = 177642.78
$ws.Range("M132").Value = -175112.78

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1171.0834
$ws.Range("I16").Value = 1205.8485
$ws.Range("K16").Value = 1205.8485
$ws.Range("M16").Value = -1035.8485
$ws.Range("H46").Value = 669882.8
$ws.Range("I46").Value = 1668625
$ws.Range("K46").Value = 1668625
$ws.Range("M46").Value = -1668437

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 16231.192
$ws.Range("I132").Value = 10881.75
$ws.Range("J132").Value = 20816.428
$ws.Range("K132").Value = 32645.25
$ws.Range("L132").Value = 62449.284
$ws.Range("M132").Value = -30115.25
$ws.Range("N132").Value = -67509.284
$ws.Range("H138").Value = 345000
$ws.Range("J138").Value = 345000
$ws.Range("L138").Value = 345000
$ws.Range("N138").Value = -355280

